$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster data (rows 2-19), columns A (name), B (position), C (team)
$players = @(
    @{A="D'Angelo Russell";     B="PG";        C="Brooklyn Nets"},
    @{A="Norman Powell";        B="SG,SF";     C="LA Clippers"},
    @{A="LeBron James";         B="SF,PF";     C="Los Angeles Lakers"},
    @{A="Myles Turner";         B="C";         C="Indiana Pacers"},
    @{A="Trayce Jackson-Davis"; B="PF,C";      C="Golden State Warriors"},
    @{A="Gradey Dick";          B="SG,SF";     C="Toronto Raptors"},
    @{A="Jabari Smith Jr.";     B="PF,C";      C="Houston Rockets"},
    @{A="Jamal Murray";         B="PG,SG";     C="Denver Nuggets"},
    @{A="Coby White";           B="PG,SG";     C="Chicago Bulls"},
    @{A="Immanuel Quickley";    B="PG,SG";     C="Toronto Raptors"},
    @{A="Walker Kessler";       B="C";         C="Utah Jazz"},
    @{A="Devin Vassell";        B="SG,SF";     C="San Antonio Spurs"},
    @{A="Desmond Bane";         B="SG,SF";     C="Memphis Grizzlies"},
    @{A="Devin Booker";         B="PG,SG";     C="Phoenix Suns"},
    @{A="Jalen Brunson";        B="PG";        C="New York Knicks"},
    @{A="Kawhi Leonard";        B="SG,SF,PF";  C="LA Clippers"},
    @{A="Brandon Ingram";       B="SG,SF,PF";  C="New Orleans Pelicans"},
    @{A="Trae Young";           B="PG";        C="Atlanta Hawks"}
)

$row = 2
foreach ($p in $players) {
    $ws.Cells.Item($row, 1).Value = $p.A
    $ws.Cells.Item($row, 2).Value = $p.B
    $ws.Cells.Item($row, 3).Value = $p.C
    $row = $row + 1
}
